$d = $word.ActiveDocument
$Q  = [char]0x201C   # “
$QR = [char]0x201D   # ”

function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: not found"
        Write-Output $old
    }
    return $found
}

# ---------------------------------------------------------------------------
# Change 1: merge "...and the acceptor(s)." paragraph with the following
# "Some subdirectories Pc, Zing..." paragraph, dropping the sentence about
# testconfig*.txt file naming.
# ---------------------------------------------------------------------------
$old1 = "The names of the testconfig*.txt files for Pc, Zing and Prt are testconfigPc.txt, testconfigZing.txt and testconfigPrt.txt, respectively.^p"
$new1 = ""
Replace-Text $old1 $new1 | Out-Null

# ---------------------------------------------------------------------------
# Change 2: drop the _GoBack bookmark splitting the SampleProtocols.txt run
# (replacing the text with itself collapses any bookmark inside the range).
# ---------------------------------------------------------------------------
$old2 = $Q + "testP.bat SampleProtocols.txt" + $QR
$new2 = $Q + "testP.bat SampleProtocols.txt" + $QR
Replace-Text $old2 $new2 | Out-Null

# ---------------------------------------------------------------------------
# Change 7: remove the "*" between "testconfig" and ".txt" in the
# "defined in testconfig*.txt" sentence, and re-create the _GoBack bookmark
# at that exact point (it used to live around "SampleProtocols.txt").
# ---------------------------------------------------------------------------
$old3 = "defined in " + $Q + "testconfig*.txt" + $QR
$new3 = "defined in " + $Q + "testconfig.txt" + $QR
Replace-Text $old3 $new3 | Out-Null

$find2 = $d.Content.Find
$search2 = "defined in " + $Q + "testconfig"
$found2 = $find2.Execute($search2, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r = $find2.Parent.Duplicate
    $r.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $r) | Out-Null
} else {
    Write-Output "WARNING: could not locate insertion point for _GoBack bookmark"
}

# ---------------------------------------------------------------------------
# Change 9: drop the "(where "*" stands for ...)" parenthetical.
# ---------------------------------------------------------------------------
$old4 = "create a file called " + $Q + "testconfig*.txt" + $QR + ", which will define how to run your test case (where " + $Q + "*" + $QR + " stands for " + $Q + "Pc" + $QR + ", " + $Q + "Prt" + $QR + " or " + $Q + "Zing" + $QR + ")."
$new4 = "create a file called " + $Q + "testconfig.txt" + $QR + ", which will define how to run your test case."
Replace-Text $old4 $new4 | Out-Null

# ---------------------------------------------------------------------------
# Change 12: collapse the long per-tool directive list down to generic
# "arg" / "inc" directives, and drop the "acc" / "igp" directives.
# ---------------------------------------------------------------------------
$old5 = "runPc: pc.exe to run. Must be exactly one such directive.^p" + `
        "runZing: zinger.exe to run. Must be exactly one such directive.^p" + `
        "runPrt: runtime exe to run. Must be exactly one such directive.^p" + `
        "argPc: An arg to pass to pc.exe. If more than one arg directive, then args are passed in order^p" + `
        "argZing: An arg to pass to zinger.exe. If more than one arg directive, then args are passed in order^p" + `
        "argPrt: An arg to pass to runtime. If more than one arg directive, then args are passed in order^p" + `
        "incPc: A file that should be included as output for pc.exe. Can be more than one such directive.^p" + `
        "incZing: A file that should be included as output for zinger.exe. Can be more than one such directive.^p" + `
        "incPrt: A file that should be included as output for runtime. Can be more than one such directive.^p" + `
        "acc: A directory containing acceptor files (more about this later). Must be exactly one such directive."
$new5 = "arg: An arg to pass to pc.exe, runtime or zinger.exe. If more than one arg directive, then args are passed in order^p" + `
        "inc: A file that should be included as output for pc.exe, runtime or zinger.exe. Can be more than one such directive."
Replace-Text $old5 $new5 | Out-Null

$old6 = "igp: Ignores output sent to the prompt by run^p"
$new6 = ""
Replace-Text $old6 $new6 | Out-Null

# ---------------------------------------------------------------------------
# Change 5: fix "IN the latter case" capitalization typo.
# ---------------------------------------------------------------------------
$old7 = "directory. IN the latter case"
$new7 = "directory. In the latter case"
Replace-Text $old7 $new7 | Out-Null

# ---------------------------------------------------------------------------
# Remaining "testconfig*" -> "testconfig" blanket simplification (all
# remaining occurrences, after the special-cased ones above already went
# through their own handling).
# ---------------------------------------------------------------------------
$old8 = "testconfig*"
$new8 = "testconfig"
$find3 = $d.Content.Find
$find3.Execute($old8, $false, $false, $false, $false, $false, $true, 1, $false, $new8, 2) | Out-Null

Write-Output "done"
